$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "26.459.84"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value2 = "  +1.88%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "1.669.58"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value2 = "  +1.46%  "
$ws.Range("E4").Value2 = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "219.70"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value2 = "  +2.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "0.5276"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value2 = "  +1.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.2673"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value2 = "  +2.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.06366"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value2 = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "21.76"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value2 = "  +5.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.07813"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value2 = "  +1.54%  "
$ws.Range("B12").Value2 = "Polkadot"
$ws.Range("C12").Value2 = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "4.478"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value2 = "  +1.19%  "
$ws.Range("B13").Value2 = "WrappedEther"
$ws.Range("C13").Value2 = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "1.669.64"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value2 = "  +1.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "0.5546"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value2 = "  +0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.0₅8300"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value2 = "  +0.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "65.40"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value2 = "  +1.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "26.472.13"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value2 = "  +1.88%  "
$ws.Range("E18").Value2 = "  +0.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "4.744"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value2 = "  +0.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "192.33"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value2 = "  +1.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "10.32"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value2 = "  +1.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "6.290"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value2 = "  +0.45%  "
$ws.Range("E23").Value2 = "  +0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "0.1263"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value2 = "  +1.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "138.16"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value2 = "  -3.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "7.390"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value2 = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "16.26"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value2 = "  +2.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "1.417"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value2 = "  +1.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "0.06187"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value2 = "  +4.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "1.293"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value2 = "  +2.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "3.609"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value2 = "  +6.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "3.419"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value2 = "  +0.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "1.678"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value2 = "  +2.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "1.004"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value2 = "  +1.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "0.6049"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value2 = "  +7.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "2.417"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value2 = "  +1.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "2.772"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value2 = "  +0.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.01610"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value2 = "  +0.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "6.027"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value2 = "  +2.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "1.090.02"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value2 = "  +5.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.8553"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value2 = "  +0.05%  "
$ws.Range("E42").Value2 = "  +0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "100.65"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value2 = "  +2.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "1.813.90"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value2 = "  +1.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "58.10"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value2 = "  +4.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "0.0₈109"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value2 = "  -1.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "8.099"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value2 = "  +0.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.9960"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value2 = "  -0.47%  "
$ws.Range("B49").Value2 = "Cronos"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.05204"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value2 = "  +1.13%  "
$ws.Range("B50").Value2 = "RenderToken"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "1.470"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value2 = "  +7.27%  "
$ws.Range("B51").Value2 = "Mantle"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "0.4231"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value2 = "  +0.33%  "
